# Apply crypto price/volume updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.516.68'
$ws.Range('E2').Value = '  -0.28%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.811.08'
$ws.Range('E3').Value = '  +0.61%  '

# Row 4
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.47'
$ws.Range('E5').Value = '  +0.42%  '

# Row 7
$ws.Range('E7').Value = '  +0.11%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '34.99'
$ws.Range('E8').Value = '  +6.65%  '

# Row 9
$ws.Range('E9').Value = '  +2.37%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0697'
$ws.Range('E10').Value = '  +0.15%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0957'
$ws.Range('E11').Value = '  +0.70%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.073.39'
$ws.Range('E12').Value = '  +0.62%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.821.60'
$ws.Range('E13').Value = '  +1.09%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.22'
$ws.Range('E14').Value = '  +0.58%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.651'
$ws.Range('E15').Value = '  +1.51%  '

# Row 16
$ws.Range('E16').Value = '  +4.44%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '34.522.56'
$ws.Range('E17').Value = '  -0.18%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.28'

# Row 19
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '246.43'
$ws.Range('E19').Value = '  -0.31%  '

# Row 20
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').Value = '  -0.58%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.46'
$ws.Range('E21').Value = '  +0.67%  '

# Row 22
$ws.Range('E22').Value = '  +0.24%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.18'
$ws.Range('E23').Value = '  -0.09%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '172.70'
$ws.Range('E24').Value = '  -0.19%  '

# Row 25
$ws.Range('E25').Value = '  +2.04%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.10'
$ws.Range('E26').Value = '  +10.51%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.82'
$ws.Range('E27').Value = '  +1.12%  '

# Row 28
$ws.Range('E28').Value = '  +2.47%  '

# Row 29
$ws.Range('E29').Value = '  -0.09%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.03'
$ws.Range('E30').Value = '  -0.17%  '

# Row 31
$ws.Range('E31').Value = '  +1.91%  '

# Row 32
$ws.Range('E32').Value = '  +1.72%  '

# Row 33
$ws.Range('E33').Value = '  +0.39%  '

# Row 34
$ws.Range('E34').Value = '  +0.44%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.398.92'
$ws.Range('E35').Value = '  -2.42%  '

# Row 36
$ws.Range('E36').Value = '  +0.52%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.48'
$ws.Range('E37').Value = '  -3.37%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').Value = '  -0.06%  '

# Row 39
$ws.Range('E39').Value = '  -0.16%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '83.83'
$ws.Range('E40').Value = '  -1.21%  '

# Row 41
$ws.Range('E41').Value = '  +1.88%  '

# Row 42
$ws.Range('E42').Value = '  +2.73%  '

# Row 43
$ws.Range('E43').Value = '  -0.13%  '

# Row 44
$ws.Range('E44').Value = '  +5.72%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.35'
$ws.Range('E45').Value = '  -3.72%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0514'
$ws.Range('E46').Value = '  -2.32%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.02'
$ws.Range('E47').Value = '  -1.52%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.972.95'
$ws.Range('E48').Value = '  +0.58%  '

# Row 49
$ws.Range('E49').Value = '  +0.26%  '

# Row 50
$ws.Range('E50').Value = '  +2.04%  '

# Row 51
$ws.Range('E51').Value = '  +0.06%  '
